$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$ws.Range("B2").Value = 0.05225407675749848
$ws.Range("C2").Value = 0.3368560989598441
$ws.Range("D2").Value = 0.1992170795427886
$ws.Range("E2").Value = 0.446337405493634
$ws.Range("F2").Value = 0.4600010262589005

$ws.Range("B3").Value = 0.08987703397608017
$ws.Range("C3").Value = 0.5080841183286242
$ws.Range("D3").Value = 0.486212353313983
$ws.Range("E3").Value = 0.6972892895448652
$ws.Range("F3").Value = 0.7197075874391274

$ws.Range("B4").Value = 0.06144343887958807
$ws.Range("C4").Value = 0.515660572847691
$ws.Range("D4").Value = 0.5847517101659792
$ws.Range("E4").Value = 0.7646905976707045
$ws.Range("F4").Value = 0.796110828729013

$ws.Range("B5").Value = -0.04108891445372104
$ws.Range("C5").Value = 0.4747591863143105
$ws.Range("D5").Value = 0.3885413735943449
$ws.Range("E5").Value = 0.6233308700797233
$ws.Range("F5").Value = 0.6523330301109209

$ws.Range("B6").Value = -0.09832660749547555
$ws.Range("C6").Value = 0.4587986694775097
$ws.Range("D6").Value = 0.2851308280455349
$ws.Range("E6").Value = 0.5339764302340834
$ws.Range("F6").Value = 0.5532356402755275

$ws.Range("B7").Value = -0.06808104366449814
$ws.Range("C7").Value = 0.4672898150270331
$ws.Range("D7").Value = 0.3267507444069591
$ws.Range("E7").Value = 0.5716211546181256
$ws.Range("F7").Value = 0.6019802159440755

$ws.Range("B8").Value = -0.172913367373134
$ws.Range("C8").Value = 0.5777194860089191
$ws.Range("D8").Value = 0.4542584475935912
$ws.Range("E8").Value = 0.6739869788012163
$ws.Range("F8").Value = 0.7136044408303032

$ws.Range("B9").Value = -0.4453154589450415
$ws.Range("C9").Value = 0.5137472253348959
$ws.Range("D9").Value = 0.495705757864413
$ws.Range("E9").Value = 0.7040637455972385
$ws.Range("F9").Value = 0.6679070667641346

$ws.Range("B10").Value = -0.1135594886937667
$ws.Range("C10").Value = 0.1135594886937667
$ws.Range("D10").Value = 0.01289575747238973
$ws.Range("E10").Value = 0.1135594886937667
